# Finished Week 13 logging
# Update the "R" row (row 3) target depth stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 398
$wsOff.Range("C3").Value = 257
$wsOff.Range("D3").Value = 66
$wsOff.Range("E3").Value = 33
$wsOff.Range("G3").Value = 3

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 412
$wsDef.Range("C3").Value = 313
$wsDef.Range("D3").Value = 78
$wsDef.Range("E3").Value = 33
